$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Planifier les activités" bullet from F9 down to F12
$ws.Range("F9").Value = $null
$ws.Range("F12").Value = "▸Planifier les activités"

$ws.Range("A12").Value = "Projet Web - Intranet"
$ws.Range("B12").Value = "06/01/2025 au 18/04/2025"
